# WiscSIMSrunSTDS.xlsx - integrate uncommitted changes:
#  - drop the StdType column (E: Run/Calib/RunCalib)
#  - update the REGEX column (D) from the old "\\D*" wildcard syntax to ".*"
#  - resize columns A-D, move the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove column E (StdType) entirely - this also shrinks the used range
#    from A1:E19 down to A1:D19.
$ws.Columns("E:E").Delete()

# 2) Rewrite the REGEX strings in column D: "\\D*" -> ".*"
#    (the stored text uses two literal backslashes before the D*).
$bs = [char]92
$oldPattern = "$bs$bs" + "D*"
for ($r = 2; $r -le 19; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value2
    if ($val -ne $null) {
        $cell.Value = $val.Replace($oldPattern, ".*")
    }
}

# 3) New column widths for A:D
$ws.Columns.Item(1).ColumnWidth = 15.666666666666666
$ws.Columns.Item(2).ColumnWidth = 8.498697916666666
$ws.Columns.Item(3).ColumnWidth = 9.330729166666666
$ws.Columns.Item(4).ColumnWidth = 14.998697916666666

# 4) Move the active selection to G24
$ws.Range("G24").Select()
